$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for E..T (16 columns): ligand/receptor stats + derived edge weights
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array: row number followed by the 16 new column values (E..T) for that row,
# reflecting updated ligand/receptor TPM-derived statistics and recomputed edge weights.
$data = @(
    @(2, 2, 0.6666666666666666, 0.098365, 0.295095, 0.001274112831990382, 0.001274112831990381, 1, 0.3333333333333333, 0.041215, 0.123645, 0.01018270982542349, 0.0101827098254235, 0.004054113475, 0.036487021275, 0.00001297392125300661, 0.00001297392125300661),
    @(3, 2, 0.6666666666666666, 0.098365, 0.295095, 0.001274112831990382, 0.001274112831990381, 3, 1, 0.987332, 2.961996, 0.243933404278904, 0.2439334042789041, 0.09711891217999999, 0.87407020962, 0.000310798680542849, 0.0003107986805428491),
    @(4, 2, 0.6666666666666666, 0.098365, 0.295095, 0.001274112831990382, 0.001274112831990381, 3, 1, 0.759735, 2.279205, 0.1877025609418444, 0.1877025609418445, 0.074731333275, 0.6725819994750001, 0.0002391542414934606, 0.0002391542414934606),
    @(5, 2, 0.6666666666666666, 0.098365, 0.295095, 0.001274112831990382, 0.001274112831990381, 3, 1, 2.259265333333333, 6.777796, 0.5581813249538279, 0.558181324953828, 0.2222326345133333, 2.00009371062, 0.0007111859887010651, 0.0007111859887010651),
    @(6, 1, 0.3333333333333333, 0.03342533333333333, 0.100276, 0.0004329552799629527, 0.0004329552799629526, 1, 0.3333333333333333, 0.041215, 0.123645, 0.01018270982542349, 0.0101827098254235, 0.001377625113333334, 0.01239862602, 0.000004408657983247737, 0.000004408657983247737),
    @(7, 1, 0.3333333333333333, 0.03342533333333333, 0.100276, 0.0004329552799629527, 0.0004329552799629526, 3, 1, 0.987332, 2.961996, 0.243933404278904, 0.2439334042789041, 0.03300190121066666, 0.297017110896, 0.000105612255341889, 0.000105612255341889),
    @(8, 1, 0.3333333333333333, 0.03342533333333333, 0.100276, 0.0004329552799629527, 0.0004329552799629526, 3, 1, 0.759735, 2.279205, 0.1877025609418444, 0.1877025609418445, 0.02539439562, 0.22854956058, 0.00008126681482233943, 0.00008126681482233943),
    @(9, 1, 0.3333333333333333, 0.03342533333333333, 0.100276, 0.0004329552799629527, 0.0004329552799629526, 3, 1, 2.259265333333333, 6.777796, 0.5581813249538279, 0.558181324953828, 0.07551669685511112, 0.6796502716960001, 0.0002416675518154764, 0.0002416675518154764),
    @(10, 3, 1, 77.07094833333333, 231.212845, 0.9982929318880467, 0.9982929318880466, 1, 0.3333333333333333, 0.041215, 0.123645, 0.01018270982542349, 0.0101827098254235, 3.176479135558333, 28.588312220025, 0.01016532724618724, 0.01016532724618724),
    @(11, 3, 1, 77.07094833333333, 231.212845, 0.9982929318880467, 0.9982929318880466, 3, 1, 0.987332, 2.961996, 0.243933404278904, 0.2439334042789041, 76.09461355984666, 684.85152203862, 0.2435169933430193, 0.2435169933430193),
    @(12, 3, 1, 77.07094833333333, 231.212845, 0.9982929318880467, 0.9982929318880466, 3, 1, 0.759735, 2.279205, 0.1877025609418444, 0.1877025609418445, 58.55349693202501, 526.981472388225, 0.1873821398855286, 0.1873821398855286),
    @(13, 3, 1, 77.07094833333333, 231.212845, 0.9982929318880467, 0.9982929318880466, 3, 1, 2.259265333333333, 6.777796, 0.5581813249538279, 0.558181324953828, 174.1237217766244, 1567.11349598962, 0.5572284714133113, 0.5572284714133114)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$colLetter$r").Value = $value
    }
}
